$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.981.92"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "3.226.45"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.594"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.24%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "3.226.18"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.603"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.132"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "3.743.27"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "3.230.95"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "62.915.60"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.964"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "365.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "631.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.26%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.375"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0710"
$ws.Range("E41").Value = "  +11.70%  "
$ws.Range("E42").Value = "  +1.93%  "
$ws.Range("D43").Value = "2.867.22"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("E44").Value = "  +10.45%  "
$ws.Range("E45").Value = "  +7.78%  "
$ws.Range("E46").Value = "  +4.78%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.58%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0393"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.41%  "
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.36%  "
